# Add a new "TestUser" manager row (row 4) to the manager list sheet,
# matching the fields used by the other rows (Name, NRIC, Age, Marital
# Status, Password).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "TestUser"
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = 55
$ws.Range("D4").Value = "Married"
$ws.Range("E4").Value = 1

# Widen column B slightly (matches the author's manual column resize).
$ws.Columns.Item(2).ColumnWidth = 13.14

# Leave the selection where the author left it after entering the row.
$ws.Range("E6").Select() | Out-Null
